$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 3 values (previously row 4's data)
$ws.Range("D3").Value = 44692
$ws.Range("J3").Value = 120

# New row 4 values (previously row 5's data)
$ws.Range("D4").Value = 44221
$ws.Range("J4").Value = 250
$ws.Range("K4").Value = 1300
$ws.Range("L4").Value = 1500
$ws.Range("M4").Value = 1420
$ws.Range("N4").Value = "$/atado"
$ws.Range("O4").Value = "Provincia de Diguillín"
$ws.Range("P4").Value = 1420
$ws.Range("Q4").Value = 1

# New row 5 values (previously row 3's data)
$ws.Range("D5").Value = 44691
$ws.Range("J5").Value = 100
$ws.Range("K5").Value = 3000
$ws.Range("L5").Value = 3500
$ws.Range("M5").Value = 3250
$ws.Range("N5").Value = "$/docena de matas"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 542
$ws.Range("Q5").Value = 6
